$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

# Row 11 - Quoziente (division) involving B4
$ws.Range("H11").Formula = "=B1/B4"
$ws.Range("I11").Formula = "=B2/B4"
$ws.Range("J11").Formula = "=B3/B4"
$ws.Range("K11").Formula = "=B4/B1"
$ws.Range("L11").Formula = "=B4/B2"
$ws.Range("M11").Formula = "=B4/B3"

# Row 12 - Potenza (power) involving B4
$ws.Range("H12").Formula = "=B1^B4"
$ws.Range("I12").Formula = "=B2^B4"
$ws.Range("J12").Formula = "=B3^B4"
$ws.Range("K12").Formula = "=B4^B3"
$ws.Range("L12").Formula = "=B4^B2"
$ws.Range("M12").Formula = "=B4^B1"

# Row 13 - Quadrato (square) of B4
$ws.Range("E13").Formula = "=B4^2"

# Row 14 - Cubo (cube) of B4
$ws.Range("E14").Formula = "=B4^3"
